# Update cell E8 text ("Good Morning" -> "GIT UPDATE") and select that cell,
# matching the commit "update file with jgit".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
